$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = -12.1537
$ws.Range("D6").Value = -8.162999999999998
$ws.Range("D7").Value = -7.553699999999993
$ws.Range("C8").Value = -12.3402
$ws.Range("D8").Value = -8.236499999999994
$ws.Range("E11").Value = 13.273
$ws.Range("A12").Value = -21.61160000000001
$ws.Range("C12").Value = -12.5849
$ws.Range("C14").Value = -12.35509999999999
$ws.Range("E14").Value = 13.45400000000001
$ws.Range("D19").Value = -7.863799999999992
$ws.Range("E19").Value = 13.99389999999999
$ws.Range("D21").Value = -7.320499999999998
$ws.Range("E21").Value = 13.9716
$ws.Range("C22").Value = -11.30449999999999
$ws.Range("D24").Value = -8.221799999999991
